$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.893996000289917
$ws.Range("E2").Value = 254.2706696112509
$ws.Range("F2").Value = 0.00839555860285668
$ws.Range("G2").Value = 0.007317598813832849
$ws.Range("H2").Value = 0.006777562160607
$ws.Range("I2").Value = 0.006574432875336125
$ws.Range("J2").Value = 0.006048304449325841
$ws.Range("K2").Value = 0.005949167639074078
$ws.Range("L2").Value = 0.005867448342101447
$ws.Range("M2").Value = 0.005681221723576779
$ws.Range("N2").Value = 0.005443680159604796
$ws.Range("O2").Value = 0.005349538687709705
$ws.Range("P2").Value = 0.005263998167949615
$ws.Range("Q2").Value = 0.005178018518343672
$ws.Range("R2").Value = 0.005111267295363178
$ws.Range("S2").Value = 0.005068157494438424
$ws.Range("T2").Value = 0.005053482205065852
$ws.Range("U2").Value = 0.005030054674675074
$ws.Range("V2").Value = 0.004992287320865127
$ws.Range("W2").Value = 0.004978389784957105
$ws.Range("X2").Value = 0.004958482620759116
$ws.Range("Y2").Value = 0.004956543267275845

$ws.Range("C3").Value = 0.9819991588592529
$ws.Range("E3").Value = 255.7576862166443
$ws.Range("F3").Value = 0.00855709625387328
$ws.Range("G3").Value = 0.0073906371341676
$ws.Range("H3").Value = 0.006829220883890151
$ws.Range("I3").Value = 0.006417108989448499
$ws.Range("J3").Value = 0.005972092946303709
$ws.Range("K3").Value = 0.005972092946303709
$ws.Range("L3").Value = 0.005623144531174283
$ws.Range("M3").Value = 0.005538980301388322
$ws.Range("N3").Value = 0.005346766489033825
$ws.Range("O3").Value = 0.005346766489033825
$ws.Range("P3").Value = 0.005295347073378718
$ws.Range("Q3").Value = 0.00523639899976222
$ws.Range("R3").Value = 0.005182132683253726
$ws.Range("S3").Value = 0.005166279296265544
$ws.Range("T3").Value = 0.00510208695496064
$ws.Range("U3").Value = 0.005077811508431552
$ws.Range("V3").Value = 0.005030977237723147
$ws.Range("W3").Value = 0.00501340099789688
$ws.Range("X3").Value = 0.004993724607915987
$ws.Range("Y3").Value = 0.004985529945743552

$ws.Range("C4").Value = 0.8830037117004395
$ws.Range("E4").Value = 258.0011464708459
$ws.Range("F4").Value = 0.008750450504240221
$ws.Range("G4").Value = 0.007428825686975189
$ws.Range("H4").Value = 0.006832313793962937
$ws.Range("I4").Value = 0.006581528988632344
$ws.Range("J4").Value = 0.00628385204708888
$ws.Range("K4").Value = 0.006092395188772027
$ws.Range("L4").Value = 0.005686891702940346
$ws.Range("M4").Value = 0.005653695681661591
$ws.Range("N4").Value = 0.005456957939559233
$ws.Range("O4").Value = 0.005283720542617467
$ws.Range("P4").Value = 0.005278995103781375
$ws.Range("Q4").Value = 0.005214761884832879
$ws.Range("R4").Value = 0.005171031595024177
$ws.Range("S4").Value = 0.005119032541116684
$ws.Range("T4").Value = 0.005119032541116684
$ws.Range("U4").Value = 0.005101553902029912
$ws.Range("V4").Value = 0.005061313413051016
$ws.Range("W4").Value = 0.005060682879978267
$ws.Range("X4").Value = 0.005046934722185288
$ws.Range("Y4").Value = 0.00502926211444144

$ws.Range("C5").Value = 0.9000015258789062
$ws.Range("E5").Value = 262.5502477866321
$ws.Range("F5").Value = 0.008539685116412517
$ws.Range("G5").Value = 0.007526497218599611
$ws.Range("H5").Value = 0.006971269499497626
$ws.Range("I5").Value = 0.0064905539500391
$ws.Range("J5").Value = 0.006379847505792928
$ws.Range("K5").Value = 0.006094395720820083
$ws.Range("L5").Value = 0.005894024909914334
$ws.Range("M5").Value = 0.0058332505205918
$ws.Range("N5").Value = 0.005689782816307005
$ws.Range("O5").Value = 0.005580906141016486
$ws.Range("P5").Value = 0.005468711567510788
$ws.Range("Q5").Value = 0.005426260864235471
$ws.Range("R5").Value = 0.005347799323207171
$ws.Range("S5").Value = 0.005286937088665432
$ws.Range("T5").Value = 0.005247448089684994
$ws.Range("U5").Value = 0.005196631731403676
$ws.Range("V5").Value = 0.005157006791149944
$ws.Range("W5").Value = 0.005140726356064696
$ws.Range("X5").Value = 0.005131480039888957
$ws.Range("Y5").Value = 0.005117938553345654

$ws.Range("C6").Value = 0.8009951114654541
$ws.Range("E6").Value = 262.4470893144517
$ws.Range("F6").Value = 0.008583839968238168
$ws.Range("G6").Value = 0.007464981996728706
$ws.Range("H6").Value = 0.006951044305928853
$ws.Range("I6").Value = 0.00650527101108664
$ws.Range("J6").Value = 0.006365203384811657
$ws.Range("K6").Value = 0.006067305846395559
$ws.Range("L6").Value = 0.005852948093393954
$ws.Range("M6").Value = 0.005818274477339521
$ws.Range("N6").Value = 0.005490828533731356
$ws.Range("O6").Value = 0.005465684042434612
$ws.Range("P6").Value = 0.005441374846924273
$ws.Range("Q6").Value = 0.005384306316691701
$ws.Range("R6").Value = 0.005349044636173719
$ws.Range("S6").Value = 0.005317118558449183
$ws.Range("T6").Value = 0.005263399757991952
$ws.Range("U6").Value = 0.005203448491560918
$ws.Range("V6").Value = 0.005186793674804797
$ws.Range("W6").Value = 0.005150076047159326
$ws.Range("X6").Value = 0.005137111824741557
$ws.Range("Y6").Value = 0.005115927666948374

$ws.Range("C7").Value = 0.943011999130249
$ws.Range("E7").Value = 262.9002613004322
$ws.Range("F7").Value = 0.008768216552876662
$ws.Range("G7").Value = 0.007620220398347506
$ws.Range("H7").Value = 0.007012680334410742
$ws.Range("I7").Value = 0.006546069182156401
$ws.Range("J7").Value = 0.006225098410961602
$ws.Range("K7").Value = 0.006046079858083527
$ws.Range("L7").Value = 0.005812433871661094
$ws.Range("M7").Value = 0.005615278097406752
$ws.Range("N7").Value = 0.005615278097406752
$ws.Range("O7").Value = 0.005537904797729317
$ws.Range("P7").Value = 0.005410690937746566
$ws.Range("Q7").Value = 0.005379151999298656
$ws.Range("R7").Value = 0.005349341345135344
$ws.Range("S7").Value = 0.005293971928667159
$ws.Range("T7").Value = 0.005256864578727594
$ws.Range("U7").Value = 0.005220410840461625
$ws.Range("V7").Value = 0.005170181691043117
$ws.Range("W7").Value = 0.005150118979842965
$ws.Range("X7").Value = 0.005139468539635497
$ws.Range("Y7").Value = 0.005124761428858328

$ws.Range("C8").Value = 0.8500008583068848
$ws.Range("E8").Value = 257.7230928099325
$ws.Range("F8").Value = 0.008473533276012753
$ws.Range("G8").Value = 0.007558127499091636
$ws.Range("H8").Value = 0.006863469102840378
$ws.Range("I8").Value = 0.006284846350559713
$ws.Range("J8").Value = 0.006130566890047705
$ws.Range("K8").Value = 0.005879239545718493
$ws.Range("L8").Value = 0.005745509062340524
$ws.Range("M8").Value = 0.005619181310348739
$ws.Range("N8").Value = 0.005494316018819332
$ws.Range("O8").Value = 0.005401938260586187
$ws.Range("P8").Value = 0.005287154475923851
$ws.Range("Q8").Value = 0.005271121686345311
$ws.Range("R8").Value = 0.005171620337912766
$ws.Range("S8").Value = 0.005144512964393868
$ws.Range("T8").Value = 0.005101105762430847
$ws.Range("U8").Value = 0.005071870373958198
$ws.Range("V8").Value = 0.005053439084391501
$ws.Range("W8").Value = 0.005042695067480515
$ws.Range("X8").Value = 0.005030876734604835
$ws.Range("Y8").Value = 0.005023841965105896

$ws.Range("C9").Value = 0.9649703502655029
$ws.Range("E9").Value = 260.4675182507253
$ws.Range("F9").Value = 0.008578781988416892
$ws.Range("G9").Value = 0.007412699781616395
$ws.Range("H9").Value = 0.006718006907143553
$ws.Range("I9").Value = 0.00647518668510253
$ws.Range("J9").Value = 0.006069050287920372
$ws.Range("K9").Value = 0.005984949621312362
$ws.Range("L9").Value = 0.005868594454970331
$ws.Range("M9").Value = 0.005712811299334304
$ws.Range("N9").Value = 0.005599484598568624
$ws.Range("O9").Value = 0.005499210951538362
$ws.Range("P9").Value = 0.005352787694859177
$ws.Range("Q9").Value = 0.005308281785186103
$ws.Range("R9").Value = 0.005236070870734982
$ws.Range("S9").Value = 0.005217604046448476
$ws.Range("T9").Value = 0.005166546349187308
$ws.Range("U9").Value = 0.005155207867951121
$ws.Range("V9").Value = 0.005132138797035066
$ws.Range("W9").Value = 0.005106033409024144
$ws.Range("X9").Value = 0.005087862434947308
$ws.Range("Y9").Value = 0.005077339537051175

$ws.Range("C10").Value = 0.7440032958984375
$ws.Range("E10").Value = 260.866412852236
$ws.Range("F10").Value = 0.008451024298259897
$ws.Range("G10").Value = 0.007140013801404936
$ws.Range("H10").Value = 0.006723629834367469
$ws.Range("I10").Value = 0.006463907677504347
$ws.Range("J10").Value = 0.006315288550667298
$ws.Range("K10").Value = 0.006069243358815724
$ws.Range("L10").Value = 0.00583053646154507
$ws.Range("M10").Value = 0.00583053646154507
$ws.Range("N10").Value = 0.00571005429228627
$ws.Range("O10").Value = 0.005539613368064573
$ws.Range("P10").Value = 0.005459766284222768
$ws.Range("Q10").Value = 0.005338489777168539
$ws.Range("R10").Value = 0.00528141527916642
$ws.Range("S10").Value = 0.005181056331172749
$ws.Range("T10").Value = 0.00517116456612896
$ws.Range("U10").Value = 0.005142950326080678
$ws.Range("V10").Value = 0.005098984472938097
$ws.Range("W10").Value = 0.005091671242786409
$ws.Range("X10").Value = 0.005091671242786409
$ws.Range("Y10").Value = 0.005085115260277503

$ws.Range("C11").Value = 0.769000768661499
$ws.Range("E11").Value = 255.4996603022119
$ws.Range("F11").Value = 0.008492893888960891
$ws.Range("G11").Value = 0.007373282264214398
$ws.Range("H11").Value = 0.006781786482252358
$ws.Range("I11").Value = 0.006302856045697241
$ws.Range("J11").Value = 0.005962142091477295
$ws.Range("K11").Value = 0.005834032027287958
$ws.Range("L11").Value = 0.005821103351618094
$ws.Range("M11").Value = 0.005781003763567156
$ws.Range("N11").Value = 0.005606310836693298
$ws.Range("O11").Value = 0.005606310836693298
$ws.Range("P11").Value = 0.005402132931980636
$ws.Range("Q11").Value = 0.005281346331845362
$ws.Range("R11").Value = 0.005214440307362854
$ws.Range("S11").Value = 0.00511923710532463
$ws.Range("T11").Value = 0.005104226478408774
$ws.Range("U11").Value = 0.00507014158604926
$ws.Range("V11").Value = 0.005030417194125069
$ws.Range("W11").Value = 0.005002812791706453
$ws.Range("X11").Value = 0.004990501457215717
$ws.Range("Y11").Value = 0.004980500200822842
